# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) worksheet gains three new trailing columns:
#   H: date             -> "2011-11-23"
#   I: legislator_name  -> "賴士葆"
#   J: legislator_id    -> 866
# for the header row and every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# --- Header row (row 1) ---
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data rows ---
for ($r = 2; $r -le $lastRow; $r++) {
    # Leading apostrophe forces the date-looking string to stay plain text
    # instead of being auto-parsed into an Excel date serial number.
    $ws.Cells.Item($r, 8).Value = "'2011-11-23"
    $ws.Cells.Item($r, 9).Value = "賴士葆"
    $ws.Cells.Item($r, 10).Value = 866
}

# --- Formatting: mirror the existing header / data-row styles onto the new cells ---
$ws.Range("B1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range(("H2:J" + $lastRow)).PasteSpecial(-4122)

$excel.CutCopyMode = $false
